$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the unit label in F1 ("Unit of the value" -> "Unit")
$ws.Range("F1").Value = "Unit"

# Move the active selection to G5 (matches the selection state saved in the file)
$ws.Range("G5").Select()
